# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Comercializadora del Agro de Limarí - Tuna)
# above the existing row 95, shifting the previous rows 95-98 down to 97-100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 95 (existing rows 95-98 shift down to 97-100).
$ws.Rows("95:96").Insert()

# --- New row 95: Tuna, Primera, 2022-06-02 ---
$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44714
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100107
$ws.Cells.Item(95, 8).Value = "Otros"
$ws.Cells.Item(95, 9).Value = 100107011
$ws.Cells.Item(95, 10).Value = "Tuna"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 160
$ws.Cells.Item(95, 14).Value = 14000
$ws.Cells.Item(95, 15).Value = 15000
$ws.Cells.Item(95, 16).Value = 14500
$ws.Cells.Item(95, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(95, 19).Value = 806
$ws.Cells.Item(95, 20).Value = 18

# --- New row 96: Tuna, Segunda, 2022-06-02 ---
$ws.Cells.Item(96, 1).Value = 2
$ws.Cells.Item(96, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(96, 3).Value = "Coquimbo"
$ws.Cells.Item(96, 4).Value = 44714
$ws.Cells.Item(96, 5).Value = 4
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100107
$ws.Cells.Item(96, 8).Value = "Otros"
$ws.Cells.Item(96, 9).Value = 100107011
$ws.Cells.Item(96, 10).Value = "Tuna"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 200
$ws.Cells.Item(96, 14).Value = 9000
$ws.Cells.Item(96, 15).Value = 10000
$ws.Cells.Item(96, 16).Value = 9500
$ws.Cells.Item(96, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(96, 19).Value = 528
$ws.Cells.Item(96, 20).Value = 18
